$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62 (Leve Item ID 27781)
$ws.Cells.Item(62, 8).Value = 941.5
$ws.Cells.Item(62, 9).Value = 834.1875
$ws.Cells.Item(62, 10).Value = 1800
$ws.Cells.Item(62, 11).Value = 834.1875
$ws.Cells.Item(62, 12).Value = 1800
$ws.Cells.Item(62, 13).Value = -210.1875
$ws.Cells.Item(62, 14).Value = -3048

# Row 65 (Leve Item ID 27781)
$ws.Cells.Item(65, 8).Value = 941.5
$ws.Cells.Item(65, 9).Value = 834.1875
$ws.Cells.Item(65, 10).Value = 1800
$ws.Cells.Item(65, 11).Value = 4170.9375
$ws.Cells.Item(65, 12).Value = 9000
$ws.Cells.Item(65, 13).Value = -1050.9375
$ws.Cells.Item(65, 14).Value = -15240

# Row 113 (Leve Item ID 27775)
$ws.Cells.Item(113, 8).Value = 2385
$ws.Cells.Item(113, 9).Value = 2433.3333
$ws.Cells.Item(113, 10).Value = 2240
$ws.Cells.Item(113, 11).Value = 2433.3333
$ws.Cells.Item(113, 12).Value = 2240
$ws.Cells.Item(113, 13).Value = 820.6667000000002
$ws.Cells.Item(113, 14).Value = -8748

# Row 137 (Leve Item ID 44013)
$ws.Cells.Item(137, 8).Value = 1227.2285
$ws.Cells.Item(137, 9).Value = 1031.6875
$ws.Cells.Item(137, 10).Value = 3313
$ws.Cells.Item(137, 11).Value = 3095.0625
$ws.Cells.Item(137, 12).Value = 9939
$ws.Cells.Item(137, 13).Value = -545.0625
$ws.Cells.Item(137, 14).Value = -15039

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 5168.5933
$ws.Cells.Item(32, 9).Value = 4553.1704
$ws.Cells.Item(32, 11).Value = 4553.1704
$ws.Cells.Item(32, 13).Value = -4266.1704

# Row 45 (Leve Item ID 27714)
$ws.Cells.Item(45, 8).Value = 11333.3
$ws.Cells.Item(45, 9).Value = 13792.125
$ws.Cells.Item(45, 11).Value = 13792.125
$ws.Cells.Item(45, 13).Value = -13415.125

# Row 61 (Leve Item ID 43999)
$ws.Cells.Item(61, 8).Value = 6417.522
$ws.Cells.Item(61, 9).Value = 7105.55
$ws.Cells.Item(61, 11).Value = 7105.55
$ws.Cells.Item(61, 13).Value = -6893.55

# Row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, 8).Value = 2402
$ws.Cells.Item(74, 9).Value = 2202.923
$ws.Cells.Item(74, 11).Value = 2202.923
$ws.Cells.Item(74, 13).Value = -1328.923

# Row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, 8).Value = 2402
$ws.Cells.Item(77, 9).Value = 2202.923
$ws.Cells.Item(77, 11).Value = 11014.615
$ws.Cells.Item(77, 13).Value = -6646.614999999998

# Row 110 (Leve Item ID 27708)
$ws.Cells.Item(110, 8).Value = 942.6070999999999
$ws.Cells.Item(110, 9).Value = 911.6667
$ws.Cells.Item(110, 10).Value = 1128.25
$ws.Cells.Item(110, 11).Value = 911.6667
$ws.Cells.Item(110, 12).Value = 1128.25
$ws.Cells.Item(110, 13).Value = 1133.3333
$ws.Cells.Item(110, 14).Value = -5218.25

# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 3666204.2
$ws.Cells.Item(122, 9).Value = 4275902.5
$ws.Cells.Item(122, 10).Value = 8014
$ws.Cells.Item(122, 11).Value = 12827707.5
$ws.Cells.Item(122, 12).Value = 24042
$ws.Cells.Item(122, 13).Value = -12825257.5
$ws.Cells.Item(122, 14).Value = -28942

# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 3816.1155
$ws.Cells.Item(132, 9).Value = 1709.8462
$ws.Cells.Item(132, 10).Value = 5922.385
$ws.Cells.Item(132, 11).Value = 5129.5386
$ws.Cells.Item(132, 12).Value = 17767.155
$ws.Cells.Item(132, 13).Value = -2599.5386
$ws.Cells.Item(132, 14).Value = -22827.155

# Row 136 (Leve Item ID 43999)
$ws.Cells.Item(136, 8).Value = 6417.522
$ws.Cells.Item(136, 9).Value = 7105.55
$ws.Cells.Item(136, 11).Value = 21316.65
$ws.Cells.Item(136, 13).Value = -18766.65

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Cells.Item(99, 8).Value = 250001500
$ws.Cells.Item(99, 9).Value = 333334340
$ws.Cells.Item(99, 11).Value = 333334340
$ws.Cells.Item(99, 13).Value = -333332842

# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 3912.5088
$ws.Cells.Item(134, 9).Value = 5079.4688
$ws.Cells.Item(134, 10).Value = 2418.8
$ws.Cells.Item(134, 11).Value = 15238.4064
$ws.Cells.Item(134, 12).Value = 7256.400000000001
$ws.Cells.Item(134, 13).Value = -12703.4064
$ws.Cells.Item(134, 14).Value = -12326.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 5566.8687
$ws.Cells.Item(31, 9).Value = 1400.5
$ws.Cells.Item(31, 10).Value = 12709.214
$ws.Cells.Item(31, 11).Value = 1400.5
$ws.Cells.Item(31, 12).Value = 12709.214
$ws.Cells.Item(31, 13).Value = -1105.5
$ws.Cells.Item(31, 14).Value = -13299.214

# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 5566.8687
$ws.Cells.Item(34, 9).Value = 1400.5
$ws.Cells.Item(34, 10).Value = 12709.214
$ws.Cells.Item(34, 11).Value = 1400.5
$ws.Cells.Item(34, 12).Value = 12709.214
$ws.Cells.Item(34, 13).Value = -1198.5
$ws.Cells.Item(34, 14).Value = -13113.214

# Row 58 (Leve Item ID 44021)
$ws.Cells.Item(58, 8).Value = 1629.7931
$ws.Cells.Item(58, 9).Value = 1165.625
$ws.Cells.Item(58, 11).Value = 1165.625
$ws.Cells.Item(58, 13).Value = -962.625

# Row 86 (Leve Item ID 12584)
$ws.Cells.Item(86, 8).Value = 2524.8
$ws.Cells.Item(86, 9).Value = 2635.182
$ws.Cells.Item(86, 10).Value = 2221.25
$ws.Cells.Item(86, 11).Value = 2635.182
$ws.Cells.Item(86, 12).Value = 2221.25
$ws.Cells.Item(86, 13).Value = -1512.182
$ws.Cells.Item(86, 14).Value = -4467.25

# Row 89 (Leve Item ID 12584)
$ws.Cells.Item(89, 8).Value = 2524.8
$ws.Cells.Item(89, 9).Value = 2635.182
$ws.Cells.Item(89, 10).Value = 2221.25
$ws.Cells.Item(89, 11).Value = 13175.91
$ws.Cells.Item(89, 12).Value = 11106.25
$ws.Cells.Item(89, 13).Value = -7559.91
$ws.Cells.Item(89, 14).Value = -22338.25

# Row 132 (Leve Item ID 44019)
$ws.Cells.Item(132, 8).Value = 2163.8667
$ws.Cells.Item(132, 9).Value = 1872.7084
$ws.Cells.Item(132, 10).Value = 3328.5
$ws.Cells.Item(132, 11).Value = 5618.1252
$ws.Cells.Item(132, 12).Value = 9985.5
$ws.Cells.Item(132, 13).Value = -3088.1252
$ws.Cells.Item(132, 14).Value = -15045.5

# Row 134 (Leve Item ID 44020)
$ws.Cells.Item(134, 8).Value = 5683.5386
$ws.Cells.Item(134, 9).Value = 6949.1
$ws.Cells.Item(134, 10).Value = 1465
$ws.Cells.Item(134, 11).Value = 20847.3
$ws.Cells.Item(134, 12).Value = 4395
$ws.Cells.Item(134, 13).Value = -18312.3
$ws.Cells.Item(134, 14).Value = -9465

# Row 136 (Leve Item ID 44021)
$ws.Cells.Item(136, 8).Value = 1629.7931
$ws.Cells.Item(136, 9).Value = 1165.625
$ws.Cells.Item(136, 11).Value = 3496.875
$ws.Cells.Item(136, 13).Value = -946.875

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Leve Item ID 36169)
$ws.Cells.Item(102, 8).Value = 1393.3334
$ws.Cells.Item(102, 9).Value = 1000
$ws.Cells.Item(102, 10).Value = 1590
$ws.Cells.Item(102, 11).Value = 1000
$ws.Cells.Item(102, 12).Value = 1590
$ws.Cells.Item(102, 13).Value = 622
$ws.Cells.Item(102, 14).Value = -4834

# Row 132 (Leve Item ID 44008)
$ws.Cells.Item(132, 8).Value = 4121.0386
$ws.Cells.Item(132, 9).Value = 5178.5386
$ws.Cells.Item(132, 10).Value = 3063.5386
$ws.Cells.Item(132, 11).Value = 15535.6158
$ws.Cells.Item(132, 12).Value = 9190.6158
$ws.Cells.Item(132, 13).Value = -13005.6158
$ws.Cells.Item(132, 14).Value = -14250.6158

$ws = $wb.Worksheets.Item("LTW")
# Row 13 (Leve Item ID 3546)
$ws.Cells.Item(13, 8).Value = 2399.5
$ws.Cells.Item(13, 10).Value = 2399.5
$ws.Cells.Item(13, 12).Value = 2399.5
$ws.Cells.Item(13, 14).Value = -2679.5

# Row 109 (Leve Item ID 27209)
$ws.Cells.Item(109, 8).Value = 46000
$ws.Cells.Item(109, 10).Value = 46000
$ws.Cells.Item(109, 12).Value = 46000
$ws.Cells.Item(109, 14).Value = -48774

# Row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, 8).Value = 7145557.5
$ws.Cells.Item(122, 9).Value = 23811522
$ws.Cells.Item(122, 11).Value = 71434566
$ws.Cells.Item(122, 13).Value = -71432116

# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 19104596
$ws.Cells.Item(132, 9).Value = 28655034
$ws.Cells.Item(132, 10).Value = 3724.125
$ws.Cells.Item(132, 11).Value = 85965102
$ws.Cells.Item(132, 12).Value = 11172.375
$ws.Cells.Item(132, 13).Value = -85962572
$ws.Cells.Item(132, 14).Value = -16232.375

# Row 136 (Leve Item ID 44060)
$ws.Cells.Item(136, 8).Value = 5654.6
$ws.Cells.Item(136, 9).Value = 6303.926
$ws.Cells.Item(136, 11).Value = 18911.778
$ws.Cells.Item(136, 13).Value = -16361.778

$ws = $wb.Worksheets.Item("WVR")
# Row 14 (Leve Item ID 2658)
$ws.Cells.Item(14, 8).Value = 70000
$ws.Cells.Item(14, 9).Value = 70000
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 70000
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -69832
$ws.Cells.Item(14, 14).ClearContents()

# Row 122 (Leve Item ID 36208)
$ws.Cells.Item(122, 8).Value = 2813.8262
$ws.Cells.Item(122, 9).Value = 2288.8125
$ws.Cells.Item(122, 11).Value = 6866.4375
$ws.Cells.Item(122, 13).Value = -4416.4375

# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 2332.238
$ws.Cells.Item(132, 9).Value = 1983.625
$ws.Cells.Item(132, 10).Value = 2546.7693
$ws.Cells.Item(132, 11).Value = 5950.875
$ws.Cells.Item(132, 12).Value = 7640.3079
$ws.Cells.Item(132, 13).Value = -3420.875
$ws.Cells.Item(132, 14).Value = -12700.3079

# Row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 2387.3513
$ws.Cells.Item(136, 9).Value = 2791.1738
$ws.Cells.Item(136, 10).Value = 1723.9286
$ws.Cells.Item(136, 11).Value = 8373.5214
$ws.Cells.Item(136, 12).Value = 5171.7858
$ws.Cells.Item(136, 13).Value = -5823.5214
$ws.Cells.Item(136, 14).Value = -10271.7858
